$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 5
$ws.Cells.Item(2, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(2, 3).Value2 = 0.9340000152587891
$ws.Cells.Item(2, 4).Value2 = 0.001999974250793457
$ws.Cells.Item(2, 5).Value2 = 1

$ws.Cells.Item(3, 1).Value2 = 23
$ws.Cells.Item(3, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(3, 3).Value2 = 0.9305000007152557
$ws.Cells.Item(3, 4).Value2 = 0.0004999935626983643
$ws.Cells.Item(3, 5).Value2 = 2

$ws.Cells.Item(4, 1).Value2 = 21
$ws.Cells.Item(4, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(4, 3).Value2 = 0.9305000007152557
$ws.Cells.Item(4, 4).Value2 = 0.0004999935626983643
$ws.Cells.Item(4, 5).Value2 = 2

$ws.Cells.Item(5, 1).Value2 = 31
$ws.Cells.Item(5, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(5, 3).Value2 = 0.9294999837875366
$ws.Cells.Item(5, 4).Value2 = 0.005499958992004395
$ws.Cells.Item(5, 5).Value2 = 4

$ws.Cells.Item(6, 1).Value2 = 22
$ws.Cells.Item(6, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(6, 3).Value2 = 0.9289999902248383
$ws.Cells.Item(6, 4).Value2 = 0.002999991178512573
$ws.Cells.Item(6, 5).Value2 = 5

$ws.Cells.Item(7, 1).Value2 = 30
$ws.Cells.Item(7, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(7, 3).Value2 = 0.9284999668598175
$ws.Cells.Item(7, 4).Value2 = 0.001499980688095093
$ws.Cells.Item(7, 5).Value2 = 6

$ws.Cells.Item(8, 1).Value2 = 14
$ws.Cells.Item(8, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(8, 3).Value2 = 0.9280000329017639
$ws.Cells.Item(8, 4).Value2 = 0.0009999871253967285
$ws.Cells.Item(8, 5).Value2 = 7

$ws.Cells.Item(9, 1).Value2 = 15
$ws.Cells.Item(9, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(9, 3).Value2 = 0.9280000329017639
$ws.Cells.Item(9, 4).Value2 = 0.002999961376190186
$ws.Cells.Item(9, 5).Value2 = 7

$ws.Cells.Item(10, 1).Value2 = 13
$ws.Cells.Item(10, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(10, 3).Value2 = 0.9259999990463257
$ws.Cells.Item(10, 4).Value2 = 0.0009999871253967285
$ws.Cells.Item(10, 5).Value2 = 9

$ws.Cells.Item(11, 1).Value2 = 20
$ws.Cells.Item(11, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(11, 3).Value2 = 0.925000011920929
$ws.Cells.Item(11, 4).Value2 = 0.004000008106231689
$ws.Cells.Item(11, 5).Value2 = 10

$ws.Cells.Item(12, 1).Value2 = 28
$ws.Cells.Item(12, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(12, 3).Value2 = 0.9234999120235443
$ws.Cells.Item(12, 4).Value2 = 0.001499921083450317
$ws.Cells.Item(12, 5).Value2 = 11

$ws.Cells.Item(13, 1).Value2 = 4
$ws.Cells.Item(13, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(13, 3).Value2 = 0.9230000078678131
$ws.Cells.Item(13, 4).Value2 = 0.00899997353553772
$ws.Cells.Item(13, 5).Value2 = 12

$ws.Cells.Item(14, 1).Value2 = 7
$ws.Cells.Item(14, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(14, 3).Value2 = 0.9219999611377716
$ws.Cells.Item(14, 4).Value2 = 0.001000076532363892
$ws.Cells.Item(14, 5).Value2 = 13

$ws.Cells.Item(15, 1).Value2 = 26
$ws.Cells.Item(15, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(15, 3).Value2 = 0.9214999675750732
$ws.Cells.Item(15, 4).Value2 = 0.002500057220458984
$ws.Cells.Item(15, 5).Value2 = 14

$ws.Cells.Item(16, 1).Value2 = 29
$ws.Cells.Item(16, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(16, 3).Value2 = 0.9199999868869781
$ws.Cells.Item(16, 4).Value2 = 0.01300004124641418
$ws.Cells.Item(16, 5).Value2 = 15

$ws.Cells.Item(17, 1).Value2 = 27
$ws.Cells.Item(17, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(17, 3).Value2 = 0.9160000085830688
$ws.Cells.Item(17, 4).Value2 = 0.0009999871253967285
$ws.Cells.Item(17, 5).Value2 = 16

$ws.Cells.Item(18, 1).Value2 = 24
$ws.Cells.Item(18, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(18, 3).Value2 = 0.9154999852180481
$ws.Cells.Item(18, 4).Value2 = 0.01649999618530273
$ws.Cells.Item(18, 5).Value2 = 17

$ws.Cells.Item(19, 1).Value2 = 25
$ws.Cells.Item(19, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(19, 3).Value2 = 0.913500040769577
$ws.Cells.Item(19, 4).Value2 = 0.002499967813491821
$ws.Cells.Item(19, 5).Value2 = 18

$ws.Cells.Item(20, 1).Value2 = 17
$ws.Cells.Item(20, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(20, 3).Value2 = 0.9129999876022339
$ws.Cells.Item(20, 4).Value2 = 0.008000016212463379
$ws.Cells.Item(20, 5).Value2 = 19

$ws.Cells.Item(21, 1).Value2 = 19
$ws.Cells.Item(21, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(21, 3).Value2 = 0.9124999642372131
$ws.Cells.Item(21, 4).Value2 = 0.004500031471252441
$ws.Cells.Item(21, 5).Value2 = 20

$ws.Cells.Item(22, 1).Value2 = 6
$ws.Cells.Item(22, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(22, 3).Value2 = 0.9104999303817749
$ws.Cells.Item(22, 4).Value2 = 0.003499984741210938
$ws.Cells.Item(22, 5).Value2 = 21

$ws.Cells.Item(23, 1).Value2 = 18
$ws.Cells.Item(23, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(23, 3).Value2 = 0.899499922990799
$ws.Cells.Item(23, 4).Value2 = 0.005499988794326782
$ws.Cells.Item(23, 5).Value2 = 22

$ws.Cells.Item(24, 1).Value2 = 16
$ws.Cells.Item(24, 2).Value2 = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(24, 3).Value2 = 0.8974999189376831
$ws.Cells.Item(24, 4).Value2 = 0.01550000905990601
$ws.Cells.Item(24, 5).Value2 = 23

$ws.Cells.Item(25, 1).Value2 = 12
$ws.Cells.Item(25, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(25, 3).Value2 = 0.8930000364780426
$ws.Cells.Item(25, 4).Value2 = 0.04399999976158142
$ws.Cells.Item(25, 5).Value2 = 24

$ws.Cells.Item(26, 1).Value2 = 1
$ws.Cells.Item(26, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(26, 3).Value2 = 0.8839999735355377
$ws.Cells.Item(26, 4).Value2 = 0.002000004053115845
$ws.Cells.Item(26, 5).Value2 = 25

$ws.Cells.Item(27, 1).Value2 = 8
$ws.Cells.Item(27, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(27, 3).Value2 = 0.8824999928474426
$ws.Cells.Item(27, 4).Value2 = 0.01150000095367432
$ws.Cells.Item(27, 5).Value2 = 26

$ws.Cells.Item(28, 1).Value2 = 9
$ws.Cells.Item(28, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(28, 3).Value2 = 0.8750000596046448
$ws.Cells.Item(28, 4).Value2 = 0.02100002765655518
$ws.Cells.Item(28, 5).Value2 = 27

$ws.Cells.Item(29, 1).Value2 = 3
$ws.Cells.Item(29, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(29, 3).Value2 = 0.8715000152587891
$ws.Cells.Item(29, 4).Value2 = 0.02049994468688965
$ws.Cells.Item(29, 5).Value2 = 28

$ws.Cells.Item(30, 1).Value2 = 0
$ws.Cells.Item(30, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(30, 3).Value2 = 0.8660000860691071
$ws.Cells.Item(30, 4).Value2 = 0.02800002694129944
$ws.Cells.Item(30, 5).Value2 = 29

$ws.Cells.Item(31, 1).Value2 = 11
$ws.Cells.Item(31, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(31, 3).Value2 = 0.8644999265670776
$ws.Cells.Item(31, 4).Value2 = 0.02350002527236938
$ws.Cells.Item(31, 5).Value2 = 30

$ws.Cells.Item(32, 1).Value2 = 10
$ws.Cells.Item(32, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(32, 3).Value2 = 0.8514999747276306
$ws.Cells.Item(32, 4).Value2 = 0.01650005578994751
$ws.Cells.Item(32, 5).Value2 = 31

$ws.Cells.Item(33, 1).Value2 = 2
$ws.Cells.Item(33, 2).Value2 = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(33, 3).Value2 = 0.8474999964237213
$ws.Cells.Item(33, 4).Value2 = 0.02150002121925354
$ws.Cells.Item(33, 5).Value2 = 32
